$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.62"
$ws.Range("D3").Value = "'31.57"
$ws.Range("D4").Value = "'5.102"
$ws.Range("E4").Value = "'-0.83%"
$ws.Range("D5").Value = "'0.07803"
$ws.Range("E5").Value = "'-1.72%"
$ws.Range("D6").Value = "'2.312"
$ws.Range("E6").Value = "'-11.57%"
$ws.Range("D7").Value = "'7.800"
$ws.Range("D8").Value = "'3.834"
$ws.Range("D9").Value = "'0.9144"
$ws.Range("E9").Value = "'0.62%"
$ws.Range("D10").Value = "'0.1746"
$ws.Range("E10").Value = "'0.68%"
$ws.Range("D11").Value = "'0.07511"
$ws.Range("E11").Value = "'3.61%"
$ws.Range("D12").Value = "'0.09157"
$ws.Range("E12").Value = "'13.71%"
$ws.Range("D13").Value = "'0.03102"
$ws.Range("E13").Value = "'2.57%"
$ws.Range("E14").Value = "'0.57%"
$ws.Range("D15").Value = "'0.001506"
$ws.Range("E15").Value = "'1.04%"
$ws.Range("D16").Value = "'0.005853"
$ws.Range("E16").Value = "'-2.94%"
$ws.Range("D17").Value = "'3.481"
$ws.Range("E17").Value = "'-0.66%"
$ws.Range("E20").Value = "'0.68%"
$ws.Range("D21").Value = "'4.025"
$ws.Range("E21").Value = "'-13.04%"
$ws.Range("D23").Value = "'0.04591"
$ws.Range("E23").Value = "'0.08%"
$ws.Range("D24").Value = "'0.001253"
$ws.Range("E24").Value = "'-0.55%"
$ws.Range("D25").Value = "'0.004457"
$ws.Range("E25").Value = "'0.17%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'6.12%"
$ws.Range("E27").Value = "'-1.42%"
$ws.Range("D39").Value = "'0.01774"
$ws.Range("E39").Value = "'-3.31%"
$ws.Range("D40").Value = "'0.04794"
$ws.Range("D41").Value = "'0.007412"
$ws.Range("E41").Value = "'5.52%"
$ws.Range("D42").Value = "'0.1357"
$ws.Range("E42").Value = "'1.14%"
$ws.Range("D43").Value = "'0.002192"
$ws.Range("E43").Value = "'-2.06%"
$ws.Range("D44").Value = "'0.01023"
$ws.Range("E44").Value = "'-1.91%"
$ws.Range("D45").Value = "'0.00006208"
$ws.Range("E45").Value = "'-3.21%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("E47").Value = "'28.79%"
$ws.Range("D48").Value = "'0.7452"
$ws.Range("E48").Value = "'-9.18%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.07%"
